# Weekly refresh of the fruit/vegetable price rows: the D (Fecha), J
# (Volumen), K (Precio minimo), L (Precio maximo), M (Precio promedio
# ponderado) and P (Precio $/Kg) values for rows 2, 7, 8, 9 and 10 get
# rotated to the next week's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot the current (pre-edit) values for the affected rows/columns
# before overwriting any of them.
$rows = @(2, 7, 8, 9, 10)
$cols = @("D", "J", "K", "L", "M", "P")

$snapshot = @{}
foreach ($r in $rows) {
    $snapshot[$r] = @{}
    foreach ($c in $cols) {
        $snapshot[$r][$c] = $ws.Range("$c$r").Value2
    }
}

# Cyclic reassignment observed in the diff:
#   row 2  <- old row 10
#   row 7  <- old row 8
#   row 8  <- old row 9
#   row 9  <- old row 2
#   row 10 <- old row 7
$rowSource = @{
    2  = 10
    7  = 8
    8  = 9
    9  = 2
    10 = 7
}

foreach ($r in $rows) {
    $src = $rowSource[$r]
    foreach ($c in $cols) {
        $ws.Range("$c$r").Value2 = $snapshot[$src][$c]
    }
}
